# Remove the "Branch" parameter column from the track block table.
# The "Branch" header (column H) and its data are deleted entirely;
# all columns to the right (Station, Station Side, ELEVATION (M),
# CUMALTIVE ELEVATION (M)) shift one position to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("H").Delete()

# Restore the cell selection left behind after the edit.
$ws.Range("H18").Select()
